# optimization plan.xlsx - debugs & plan update
# 1. debug： calendar改變age顯示方式時離整月不到一周時會顯示整月
# 2. debug：calendar無記錄時title顯示為1 jan 1
# 3. 增加main界面 按鍵盤回退鍵 退出程序

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Row 56/57 rework ----------------------------------------------------
# B56 ("增加輔食模塊...") moves down to C60; B56 itself is removed afterwards.
$ws.Range("B56").Cut($ws.Range("C60"))
$ws.Range("B56").Clear()

# C56 keeps its text (分辨率自適應) but picks up the "Good" style used elsewhere in column C
$ws.Range("C50").Copy()
$ws.Range("C56").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# C57 switches from the plain "Good" style to the column-C "Good" style
$ws.Range("B3").Copy()
$ws.Range("C57").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# C60 (now holding the moved text) gets the column-C "Good" style too
$ws.Range("B3").Copy()
$ws.Range("C60").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Insert two new rows before the old row 61 (old 61->63, old 66-72->68-74)
$ws.Rows("61:62").Insert()

# New C61: "debug calendar無記錄時title顯示為1 jan 1"
$ws.Range("B3").Copy()
$ws.Range("C61").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C61").Value = "debug calendar無記錄時title顯示為1 jan 1"

# New C62: "debug： calendar改變age顯示方式時離整月不到一周時會顯示整月" (bordered "Good" style)
$ws.Range("C54").Copy()
$ws.Range("C62").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C62").Value = "debug： calendar改變age顯示方式時離整月不到一周時會顯示整月"

# --- Final selection ------------------------------------------------------
$ws.Range("C62").Select()

Write-Host "edit complete"
